# Forecast summary update — shift each forecast row forward by one week
# and refresh the MyForecast numbers / derived Summary metrics.
# (Commit: "Tried to implement Penality Reward System (unfinished)")

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": Week_Start_Date (B) + MyForecast (D) ---
# Leading "'" forces Excel to keep the Week_Start_Date column as literal
# text (matching the source data) instead of auto-converting it to a date
# serial number.
$forecastRows = @(
    @{ Row = 2;  B = "2025-01-12"; D = 27 },
    @{ Row = 3;  B = "2025-01-19"; D = 28 },
    @{ Row = 4;  B = "2025-01-26"; D = 30 },
    @{ Row = 5;  B = "2025-02-02"; D = 30 },
    @{ Row = 6;  B = "2025-02-09"; D = 30 },
    @{ Row = 7;  B = "2025-02-16"; D = 29 },
    @{ Row = 8;  B = "2025-02-23"; D = 29 },
    @{ Row = 9;  B = "2025-03-02"; D = 29 },
    @{ Row = 10; B = "2025-03-09"; D = 29 },
    @{ Row = 11; B = "2025-03-16"; D = 29 },
    @{ Row = 12; B = "2025-03-23"; D = 29 },
    @{ Row = 13; B = "2025-03-30"; D = 29 },
    @{ Row = 14; B = "2025-04-06"; D = 28 },
    @{ Row = 15; B = "2025-04-13"; D = 27 },
    @{ Row = 16; B = "2025-04-20"; D = 27 },
    @{ Row = 17; B = "2025-04-27"; D = 29 }
)

foreach ($r in $forecastRows) {
    $ws1.Cells.Item($r.Row, 2).Value = "'" + $r.B
    $ws1.Cells.Item($r.Row, 4).Value = $r.D
}

# --- Sheet "Summary": refreshed aggregate metrics ---
# All Value cells on this sheet are stored as text, including the ones
# that merely look like plain numbers, so every write is apostrophe
# -escaped to keep it text (preventing Excel from re-typing them as
# numbers or dates).
$ws2.Range("B2").Value  = "2023-01-29 to 2025-01-05"
$ws2.Range("B8").Value  = "1770 units"
$ws2.Range("B9").Value  = "'459"
$ws2.Range("B10").Value = "'232"
$ws2.Range("B11").Value = "'115"
$ws2.Range("B12").Value = "'30"
$ws2.Range("B14").Value = "'27"
$ws2.Range("B15").Value = "'2025-01-12"
